# Avances al 20 de Marzo
# Fill in the missing "Ciudadano #5" (column F) answers on the "Ciudadanos" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudadanos")

# F3: repeat the same "Comprendio bien las pantallas" observation already present in E3
$ws.Range("F3").Value = $ws.Range("E3").Value()
$ws.Range("F3").WrapText = $true

# F7: mark as not answered / not applicable
$ws.Range("F7").Value = "-"
$ws.Range("F7").WrapText = $true

# E8 / F8: mark as not answered / not applicable
$ws.Range("E8").Value = "-"
$ws.Range("E8").WrapText = $true

$ws.Range("F8").Value = "-"
$ws.Range("F8").WrapText = $true

# Move the active selection to B4, as left by the author when saving
$ws.Range("B4").Select()
